$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 164, shifting existing row 164 (and below) down by one.
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row 164 with the Lyon, France entry.
$ws.Cells.Item(164, 1).Value = "LYS"
$ws.Cells.Item(164, 2).Value = "Lyon"
$ws.Cells.Item(164, 3).Value = 45.7263
$ws.Cells.Item(164, 4).Value = 5.0908
$ws.Cells.Item(164, 5).Value = "FR"
$ws.Cells.Item(164, 6).Value = "Europe"
$ws.Cells.Item(164, 7).Value = "Lyon"

# Copy the style (style index 1) from column A of a neighboring row onto the new A164 cell.
$ws.Cells.Item(163, 1).Copy()
$ws.Cells.Item(164, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
